# edit.ps1 -- apply the "Scenario" section + Success Criteria page-break-marker
# + new bullet list (numId=2) edits described by the task diff.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert the new "Scenario" block right after the "Task Instructions"
#    title paragraph (i.e. in place of the blank paragraph that currently
#    follows it), while preserving that trailing blank paragraph.
# ---------------------------------------------------------------------------

# NOTE: the four bullet paragraphs deliberately carry pStyle=ListParagraph
# but *no* <w:numPr> yet -- the numPr/numId is added in step (2) below via
# ListFormat.ApplyBulletDefault(), which also mints the backing
# abstractNum/num pair in numbering.xml. (Baking a literal numId="2" into
# this raw XML up front would make ApplyBulletDefault() treat that id as
# already taken and skip to numId="3".)
$blockXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Scenario:</w:t></w:r></w:p><w:p><w:r><w:t>Your organisation wants to introduce a new Student Learning Management System (LMS) (</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>similar to</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> Moodle or Canvas) to replace outdated paper-based and manual processes. The system should allow:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t>Online enrolment</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t>Assignment submissions</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t>Grade tracking</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t>Messaging between students and lecturers</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">You are the </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Project Manager</w:t></w:r><w:r><w:t xml:space="preserve"> and need to prepare a </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Project Charter</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

# Locate the title paragraph ("Task Instructions") and grab the blank
# paragraph immediately following it -- InsertXML replaces the *whole*
# paragraph that owns the (possibly collapsed) range, so targeting that
# blank paragraph's range swaps it for our multi-paragraph block while the
# two trailing <w:p/> markers guarantee a blank paragraph still remains
# afterwards (InsertXML folds the very last paragraph mark of the inserted
# XML into the target's own mark).
$titlePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "Task Instructions") {
        $titlePara = $p
        break
    }
}
$blankPara = $titlePara.Next()
$blankPara.Range.InsertXML($blockXml)

# ---------------------------------------------------------------------------
# 2) Apply real bullet-list numbering (pStyle=ListParagraph already set
#    above) to the four new list-item paragraphs so numId=2 is backed by an
#    actual abstractNum/num pair in numbering.xml, matching the new bullet
#    list the diff introduces.
# ---------------------------------------------------------------------------

$first = $null
$last = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd()
    if ($t -eq "Online enrolment") { $first = $p }
    if ($t -eq "Messaging between students and lecturers") { $last = $p }
}
$bulletRange = $d.Range($first.Range.Start, $last.Range.End)
$bulletRange.ListFormat.ApplyBulletDefault()

# ---------------------------------------------------------------------------
# 3) Mark the "Success Criteria" run with <w:lastRenderedPageBreak/>, as in
#    the diff, while preserving its existing bold run formatting and its
#    numbered-paragraph properties.
# ---------------------------------------------------------------------------

$successPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "Success Criteria") {
        $successPara = $p
        break
    }
}

$successXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t>Success Criteria</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$successPara.Range.InsertXML($successXml)

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
